$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (H1) onto the two new header cells so they match the rest
# of the header row's style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J for rows 2-6
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 3

$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 4
